$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 1325.8
$ws.Cells.Item(17, 10).Value = 1362
$ws.Cells.Item(17, 12).Value = 4086
$ws.Cells.Item(17, 14).Value = -4422
# Row 28
$ws.Cells.Item(28, 8).Value = 5550
$ws.Cells.Item(28, 9).Value = 4233.3335
$ws.Cells.Item(28, 11).Value = 4233.3335
$ws.Cells.Item(28, 13).Value = -3748.3335
# Row 34
$ws.Cells.Item(34, 8).Value = 5666.6665
$ws.Cells.Item(34, 9).Value = 5666.6665
$ws.Cells.Item(34, 11).Value = 5666.6665
$ws.Cells.Item(34, 13).Value = -5463.6665
# Row 36
$ws.Cells.Item(36, 8).Value = 5666.6665
$ws.Cells.Item(36, 9).Value = 5666.6665
$ws.Cells.Item(36, 11).Value = 5666.6665
$ws.Cells.Item(36, 13).Value = -4951.6665
# Row 74
$ws.Cells.Item(74, 8).Value = 7940.0625
$ws.Cells.Item(74, 9).Value = 7541.615
$ws.Cells.Item(74, 11).Value = 7541.615
$ws.Cells.Item(74, 13).Value = -6605.615
# Row 77
$ws.Cells.Item(77, 8).Value = 7940.0625
$ws.Cells.Item(77, 9).Value = 7541.615
$ws.Cells.Item(77, 11).Value = 37708.075
$ws.Cells.Item(77, 13).Value = -33028.075
# Row 116
$ws.Cells.Item(116, 8).Value = 23518.875
$ws.Cells.Item(116, 9).Value = 26069.54
$ws.Cells.Item(116, 10).Value = 12466
$ws.Cells.Item(116, 11).Value = 26069.54
$ws.Cells.Item(116, 12).Value = 12466
$ws.Cells.Item(116, 13).Value = -22627.54
$ws.Cells.Item(116, 14).Value = -19350
# Row 129
$ws.Cells.Item(129, 8).Value = 1189.5714
$ws.Cells.Item(129, 9).Value = 817.0909
$ws.Cells.Item(129, 10).Value = 2555.3333
$ws.Cells.Item(129, 11).Value = 2451.2727
$ws.Cells.Item(129, 12).Value = 7665.999899999999
$ws.Cells.Item(129, 13).Value = 2548.7273
$ws.Cells.Item(129, 14).Value = -17665.9999
# Row 132
$ws.Cells.Item(132, 8).Value = 24224.871
$ws.Cells.Item(132, 9).Value = 25564.518
$ws.Cells.Item(132, 11).Value = 76693.554
$ws.Cells.Item(132, 13).Value = -74163.554
# Row 138
$ws.Cells.Item(138, 8).Value = 30406.584
$ws.Cells.Item(138, 9).Value = 2070.9546
$ws.Cells.Item(138, 10).Value = 74934
$ws.Cells.Item(138, 11).Value = 6212.8638
$ws.Cells.Item(138, 12).Value = 224802
$ws.Cells.Item(138, 13).Value = -1072.8638
$ws.Cells.Item(138, 14).Value = -235082
# Row 141
$ws.Cells.Item(141, 8).Value = 3486.4
$ws.Cells.Item(141, 9).Value = 3811
$ws.Cells.Item(141, 10).Value = 2999.5
$ws.Cells.Item(141, 11).Value = 11433
$ws.Cells.Item(141, 12).Value = 8998.5
$ws.Cells.Item(141, 13).Value = -6253
$ws.Cells.Item(141, 14).Value = -19358.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 29112.135
$ws.Cells.Item(32, 9).Value = 33445.438
$ws.Cells.Item(32, 11).Value = 33445.438
$ws.Cells.Item(32, 13).Value = -33158.438
# Row 43
$ws.Cells.Item(43, 8).Value = 20984.3
$ws.Cells.Item(43, 10).Value = 19622
$ws.Cells.Item(43, 12).Value = 19622
$ws.Cells.Item(43, 14).Value = -20248
# Row 52
$ws.Cells.Item(52, 8).Value = 45000
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()
# Row 74
$ws.Cells.Item(74, 8).Value = 679780.9
$ws.Cells.Item(74, 9).Value = 2000670
$ws.Cells.Item(74, 10).Value = 19336.334
$ws.Cells.Item(74, 11).Value = 2000670
$ws.Cells.Item(74, 12).Value = 19336.334
$ws.Cells.Item(74, 13).Value = -1999796
$ws.Cells.Item(74, 14).Value = -21084.334
# Row 76
$ws.Cells.Item(76, 8).Value = 20630.5
$ws.Cells.Item(76, 10).Value = 21000
$ws.Cells.Item(76, 12).Value = 21000
$ws.Cells.Item(76, 14).Value = -21676
# Row 77
$ws.Cells.Item(77, 8).Value = 679780.9
$ws.Cells.Item(77, 9).Value = 2000670
$ws.Cells.Item(77, 10).Value = 19336.334
$ws.Cells.Item(77, 11).Value = 10003350
$ws.Cells.Item(77, 12).Value = 96681.67
$ws.Cells.Item(77, 13).Value = -9998982
$ws.Cells.Item(77, 14).Value = -105417.67
# Row 79
$ws.Cells.Item(79, 8).Value = 20630.5
$ws.Cells.Item(79, 10).Value = 21000
$ws.Cells.Item(79, 12).Value = 21000
$ws.Cells.Item(79, 14).Value = -23340
# Row 102
$ws.Cells.Item(102, 8).Value = 1855.8077
$ws.Cells.Item(102, 9).Value = 1466.1305
$ws.Cells.Item(102, 10).Value = 4843.3335
$ws.Cells.Item(102, 11).Value = 1466.1305
$ws.Cells.Item(102, 12).Value = 4843.3335
$ws.Cells.Item(102, 13).Value = 155.8695
$ws.Cells.Item(102, 14).Value = -8087.3335
# Row 132
$ws.Cells.Item(132, 8).Value = 1115.0454
$ws.Cells.Item(132, 9).Value = 1047.9767
$ws.Cells.Item(132, 11).Value = 3143.9301
$ws.Cells.Item(132, 13).Value = -613.9300999999996

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 1058.2
$ws.Cells.Item(99, 9).Value = 1058.2
$ws.Cells.Item(99, 11).Value = 1058.2
$ws.Cells.Item(99, 13).Value = 439.8
# Row 134
$ws.Cells.Item(134, 8).Value = 2336.52
$ws.Cells.Item(134, 9).Value = 2071.05
$ws.Cells.Item(134, 11).Value = 6213.150000000001
$ws.Cells.Item(134, 13).Value = -3678.150000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 4995.3335
$ws.Cells.Item(16, 9).Value = 4993
$ws.Cells.Item(16, 11).Value = 4993
$ws.Cells.Item(16, 13).Value = -4706
# Row 31
$ws.Cells.Item(31, 8).Value = 4546940
$ws.Cells.Item(31, 10).Value = 3039.5
$ws.Cells.Item(31, 12).Value = 3039.5
$ws.Cells.Item(31, 14).Value = -3629.5
# Row 34
$ws.Cells.Item(34, 8).Value = 4546940
$ws.Cells.Item(34, 10).Value = 3039.5
$ws.Cells.Item(34, 12).Value = 3039.5
$ws.Cells.Item(34, 14).Value = -3443.5
# Row 58
$ws.Cells.Item(58, 8).Value = 1459.9286
$ws.Cells.Item(58, 10).Value = 2136.25
$ws.Cells.Item(58, 12).Value = 2136.25
$ws.Cells.Item(58, 14).Value = -2542.25
# Row 94
$ws.Cells.Item(94, 8).Value = 3121.4
$ws.Cells.Item(94, 10).Value = 3940.5557
$ws.Cells.Item(94, 12).Value = 3940.5557
$ws.Cells.Item(94, 14).Value = -4842.5557
# Row 105
$ws.Cells.Item(105, 8).Value = 49066.75
$ws.Cells.Item(105, 9).Value = 64889
$ws.Cells.Item(105, 11).Value = 64889
$ws.Cells.Item(105, 13).Value = -63142
# Row 113
$ws.Cells.Item(113, 8).Value = 4995.3335
$ws.Cells.Item(113, 9).Value = 4993
$ws.Cells.Item(113, 11).Value = 4993
$ws.Cells.Item(113, 13).Value = -2823
# Row 132
$ws.Cells.Item(132, 8).Value = 68426.8
$ws.Cells.Item(132, 9).Value = 91809.27
$ws.Cells.Item(132, 11).Value = 275427.81
$ws.Cells.Item(132, 13).Value = -272897.81
# Row 134
$ws.Cells.Item(134, 8).Value = 3484.0645
$ws.Cells.Item(134, 9).Value = 2969.1365
$ws.Cells.Item(134, 11).Value = 8907.4095
$ws.Cells.Item(134, 13).Value = -6372.4095
# Row 136
$ws.Cells.Item(136, 8).Value = 1459.9286
$ws.Cells.Item(136, 10).Value = 2136.25
$ws.Cells.Item(136, 12).Value = 6408.75
$ws.Cells.Item(136, 14).Value = -11508.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Cells.Item(9, 8).Value = 1000
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 1000
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 3000
$ws.Cells.Item(9, 13).ClearContents()
$ws.Cells.Item(9, 14).Value = -3448
# Row 107
$ws.Cells.Item(107, 8).Value = 1617.6818
$ws.Cells.Item(107, 9).Value = 2581.5
$ws.Cells.Item(107, 10).Value = 1066.9286
$ws.Cells.Item(107, 11).Value = 7744.5
$ws.Cells.Item(107, 12).Value = 3200.7858
$ws.Cells.Item(107, 13).Value = -5824.5
$ws.Cells.Item(107, 14).Value = -7040.7858
# Row 131
$ws.Cells.Item(131, 8).Value = 117386.57
$ws.Cells.Item(131, 10).Value = 1901.8572
$ws.Cells.Item(131, 12).Value = 5705.571599999999
$ws.Cells.Item(131, 14).Value = -15785.5716
# Row 136
$ws.Cells.Item(136, 8).Value = 2300
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 74
$ws.Cells.Item(74, 8).Value = 27000
$ws.Cells.Item(74, 10).Value = 27000
$ws.Cells.Item(74, 12).Value = 27000
$ws.Cells.Item(74, 14).Value = -28872
# Row 77
$ws.Cells.Item(77, 8).Value = 27000
$ws.Cells.Item(77, 10).Value = 27000
$ws.Cells.Item(77, 12).Value = 81000
$ws.Cells.Item(77, 14).Value = -90360
# Row 113
$ws.Cells.Item(113, 8).Value = 3616.3333
$ws.Cells.Item(113, 9).Value = 3499.75
$ws.Cells.Item(113, 10).Value = 3849.5
$ws.Cells.Item(113, 11).Value = 3499.75
$ws.Cells.Item(113, 12).Value = 3849.5
$ws.Cells.Item(113, 13).Value = -1329.75
$ws.Cells.Item(113, 14).Value = -8189.5
# Row 114
$ws.Cells.Item(114, 8).Value = 77979
$ws.Cells.Item(114, 10).Value = 77979
$ws.Cells.Item(114, 12).Value = 77979
$ws.Cells.Item(114, 14).Value = -86657
# Row 122
$ws.Cells.Item(122, 8).Value = 4067.2693
$ws.Cells.Item(122, 10).Value = 4274.75
$ws.Cells.Item(122, 12).Value = 12824.25
$ws.Cells.Item(122, 14).Value = -17724.25
# Row 132
$ws.Cells.Item(132, 8).Value = 2627.353
$ws.Cells.Item(132, 9).Value = 2259.6667
$ws.Cells.Item(132, 10).Value = 3041
$ws.Cells.Item(132, 11).Value = 6779.000100000001
$ws.Cells.Item(132, 12).Value = 9123
$ws.Cells.Item(132, 13).Value = -4249.000100000001
$ws.Cells.Item(132, 14).Value = -14183

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 3909.3333
$ws.Cells.Item(40, 9).Value = 3454.8572
$ws.Cells.Item(40, 11).Value = 3454.8572
$ws.Cells.Item(40, 13).Value = -3318.8572
# Row 132
$ws.Cells.Item(132, 8).Value = 4673.3184
$ws.Cells.Item(132, 9).Value = 4590.65
$ws.Cells.Item(132, 10).Value = 5500
$ws.Cells.Item(132, 11).Value = 13771.95
$ws.Cells.Item(132, 12).Value = 16500
$ws.Cells.Item(132, 13).Value = -11241.95
$ws.Cells.Item(132, 14).Value = -21560

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 2500.75
$ws.Cells.Item(62, 10).Value = 4001.5
$ws.Cells.Item(62, 12).Value = 4001.5
$ws.Cells.Item(62, 14).Value = -5249.5
# Row 65
$ws.Cells.Item(65, 8).Value = 2500.75
$ws.Cells.Item(65, 10).Value = 4001.5
$ws.Cells.Item(65, 12).Value = 20007.5
$ws.Cells.Item(65, 14).Value = -26247.5
